# Applies the "dados na planilha fiis" change:
#  - Aportes (sheet1): add columns B..K (tipo, quantidade, preco, dy_mes,
#    dy_ano, dy_percentual, dv_ano, dv_mes, data_com, data_cadastrado) and
#    update row 2 data for BTCI11.
#  - Proventos (sheet2): update data_com date in C2.
#  - Rendimentos (sheet3): brand-new sheet summarizing the FII position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Aportes"
# ---------------------------------------------------------------------
$aportes = $wb.Worksheets.Item("Aportes")

$aportes.Range("A1").Value = "fundo"
$aportes.Range("B1").Value = "tipo"
$aportes.Range("C1").Value = "quantidade"
$aportes.Range("D1").Value = "preco"
$aportes.Range("E1").Value = "dy_mes"
$aportes.Range("F1").Value = "dy_ano"
$aportes.Range("G1").Value = "dy_percentual"
$aportes.Range("H1").Value = "dv_ano"
$aportes.Range("I1").Value = "dv_mes"
$aportes.Range("J1").Value = "data_com"
$aportes.Range("K1").Value = "data_cadastrado"

# Carry over the header styling (bold / centered / thin box border) used
# by the original header cells onto the newly inserted header cells.
$newHeaders = $aportes.Range("B1:K1")
$newHeaders.Font.Bold = $true
$newHeaders.HorizontalAlignment = -4108
$newHeaders.VerticalAlignment = -4160
$newHeaders.Borders.LineStyle = 1

$aportes.Range("A2").Value = "BTCI11"
$aportes.Range("B2").Value = "PAPEL"
$aportes.Range("C2").Value = 50
$aportes.Range("D2").Value = 9.2
$aportes.Range("E2").Value = 1.12
$aportes.Range("F2").Value = 12.27
$aportes.Range("G2").Value = 0
$aportes.Range("H2").Value = 0
$aportes.Range("I2").Value = 0.1

# "07/08/2025" is ambiguous (valid as M/D/Y) and would otherwise be
# auto-converted to a date serial; force it to stay plain text, matching
# the source data, then drop back to the default (unstyled) cell style.
$aportes.Range("J2").NumberFormat = "@"
$aportes.Range("J2").Value = "07/08/2025"
$aportes.Range("J2").Style = "Normal"

$aportes.Range("K2").Value = "22/07/2025"

# ---------------------------------------------------------------------
# Sheet "Proventos"
# ---------------------------------------------------------------------
$proventos = $wb.Worksheets.Item("Proventos")
$proventos.Range("C2").Value = "22/07/2025"

# ---------------------------------------------------------------------
# New sheet "Rendimentos" (inserted after "Proventos", i.e. at the end)
# ---------------------------------------------------------------------
$rendimentos = $wb.Worksheets.Add($null, $proventos)
$rendimentos.Name = "Rendimentos"

$rendimentos.Range("A1").Value = "FII"
$rendimentos.Range("B1").Value = "QTDE COTAS"
$rendimentos.Range("C1").Value = "ULTIMO PREÇO"
$rendimentos.Range("D1").Value = "PROVENTOS"
$rendimentos.Range("E1").Value = "RENDIMENTO MÊS"
$rendimentos.Range("F1").Value = "RENDIMENTO ANO"
$rendimentos.Range("G1").Value = "DATA COM"
$rendimentos.Range("H1").Value = "DATA QUE FOI GERADA"

$rendHeaders = $rendimentos.Range("A1:H1")
$rendHeaders.Font.Bold = $true
$rendHeaders.HorizontalAlignment = -4108
$rendHeaders.VerticalAlignment = -4160
$rendHeaders.Borders.LineStyle = 1

$rendimentos.Range("A2").Value = "BTCI11"
$rendimentos.Range("B2").Value = 50
$rendimentos.Range("C2").Value = 9.2
$rendimentos.Range("D2").Value = 0.1
$rendimentos.Range("E2").Value = 56.00000000000001
$rendimentos.Range("F2").Value = 613.5

# Same ambiguous-date fix as above for "07/08/2025".
$rendimentos.Range("G2").NumberFormat = "@"
$rendimentos.Range("G2").Value = "07/08/2025"
$rendimentos.Range("G2").Style = "Normal"

$rendimentos.Range("H2").Value = "22/07/2025"
